$d = $word.ActiveDocument

function Replace-ParagraphXml($matchText, $newXml) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text
        if ($t -eq $matchText) {
            $p.Range.InsertXML($newXml)
            return $true
        }
    }
    return $false
}

# 1) "HTML Div Element" heading -> split into tag-highlighted runs
$div_new = '<w:p w:rsidR="00A70546" w:rsidRDefault="00A70546" w:rsidP="00A70546"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t xml:space="preserve">HTML </w:t></w:r><w:r><w:t>&lt;d</w:t></w:r><w:r><w:t>iv</w:t></w:r><w:r><w:t>&gt;</w:t></w:r><w:r><w:t xml:space="preserve"> Element</w:t></w:r></w:p>'
Replace-ParagraphXml "HTML Div Element`r" $div_new

# 2) merge the img-paragraph runs that were split around the old _GoBack bookmark,
#    and drop that bookmark (it is moving to the Span heading below)
$img_new = '<w:p w:rsidR="00396668" w:rsidRDefault="00396668" w:rsidP="00396668"><w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:t>The s</w:t></w:r><w:r w:rsidRPr="00396668"><w:rPr><w:noProof/></w:rPr><w:t>ole</w:t></w:r><w:r><w:t xml:space="preserve"> &lt;img&gt; tag is not enough to display image, we need to use additional attribute src, which describes the path to the image. If the image is in the same place as the HTML document then for the path we can use just image name and extension. But if the </w:t></w:r><w:r w:rsidRPr="00396668"><w:rPr><w:noProof/></w:rPr><w:t>image</w:t></w:r><w:r><w:t xml:space="preserve"> is outside the folder where the HTML document is, we need to use this syntax:  ../image-name-and-extension.  </w:t></w:r></w:p>'
$img_match_text = "The sole <img> tag is not enough to display image, we need to use additional attribute src, which describes the path to the image. If the image is in the same place as the HTML document then for the path we can use just image name and extension. But if the image is outside the folder where the HTML document is, we need to use this syntax:  ../image-name-and-extension.  `r"
Replace-ParagraphXml $img_match_text $img_new

# 3) "HTML Span Element" heading -> split into tag-highlighted runs, carrying the
#    _GoBack bookmark that used to sit in the img paragraph above
$span_new = '<w:p w:rsidR="00910AFA" w:rsidRDefault="006131AE" w:rsidP="006131AE"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>HTML &lt;s</w:t></w:r><w:r><w:t>pan</w:t></w:r><w:r><w:t>&gt;</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> Element</w:t></w:r></w:p>'
Replace-ParagraphXml "HTML Span Element`r" $span_new
